# Raw Data from Database.xlsx - self-appraisal task rows 399-460:
# mark the relevant "Goal"/"Training" columns with a 1 for each task,
# and move the sheet's active selection to reflect where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 399
$ws.Range("J399").Value = 1

# Rows 400-409 (Goal1/Goal3/Goal5 columns)
$ws.Range("F400:F409").Value = 1
$ws.Range("H400:H409").Value = 1
$ws.Range("J400:J409").Value = 1

# Rows 410-417 (Goal1 column only)
$ws.Range("F410:F417").Value = 1

# Rows 418-427 (Goal1/Goal3/Goal5 columns)
$ws.Range("F418:F427").Value = 1
$ws.Range("H418:H427").Value = 1
$ws.Range("J418:J427").Value = 1

# Rows 428-431 (Goal4 column only)
$ws.Range("I428:I431").Value = 1

# Rows 432-436 (Goal1/Goal3 columns)
$ws.Range("F432:F436").Value = 1
$ws.Range("H432:H436").Value = 1

# Rows 437-440 (Goal1 column only)
$ws.Range("F437:F440").Value = 1

# Row 441 (Training column)
$ws.Range("K441").Value = 1

# Row 442 (Goal5 column)
$ws.Range("J442").Value = 1

# Rows 443-446 (Goal1 column only)
$ws.Range("F443:F446").Value = 1

# Rows 447-454 (Goal1 column only)
$ws.Range("F447:F454").Value = 1

# Rows 455-459 (Goal1 column only)
$ws.Range("F455:F459").Value = 1

# Row 460 (Goal4 column)
$ws.Range("I460").Value = 1

# Update the sheet's saved view state (scroll position + active selection)
$excel.ActiveWindow.ScrollRow = 457
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I464").Select()
